$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A414:A422").EntireRow.Insert()

$ws.Range("B414:D422").NumberFormat = "@"

$ws.Range("A414").Value = 1574035200
$ws.Range("B414").Value = "2019-11-18"
$ws.Range("C414").Value = "0196"
$ws.Range("D414").Value = "QES"
$ws.Range("E414").Value = 0.245
$ws.Range("F414").Value = 0.25
$ws.Range("G414").Value = 0.235
$ws.Range("H414").Value = 0.24
$ws.Range("I414").Value = 11929500

$ws.Range("A415").Value = 1574121600
$ws.Range("B415").Value = "2019-11-19"
$ws.Range("C415").Value = "0196"
$ws.Range("D415").Value = "QES"
$ws.Range("E415").Value = 0.24
$ws.Range("F415").Value = 0.24
$ws.Range("G415").Value = 0.235
$ws.Range("H415").Value = 0.235
$ws.Range("I415").Value = 3333200

$ws.Range("A416").Value = 1574208000
$ws.Range("B416").Value = "2019-11-20"
$ws.Range("C416").Value = "0196"
$ws.Range("D416").Value = "QES"
$ws.Range("E416").Value = 0.235
$ws.Range("F416").Value = 0.245
$ws.Range("G416").Value = 0.235
$ws.Range("H416").Value = 0.24
$ws.Range("I416").Value = 4628100

$ws.Range("A417").Value = 1574294400
$ws.Range("B417").Value = "2019-11-21"
$ws.Range("C417").Value = "0196"
$ws.Range("D417").Value = "QES"
$ws.Range("E417").Value = 0.24
$ws.Range("F417").Value = 0.24
$ws.Range("G417").Value = 0.235
$ws.Range("H417").Value = 0.24
$ws.Range("I417").Value = 4867600

$ws.Range("A418").Value = 1574380800
$ws.Range("B418").Value = "2019-11-22"
$ws.Range("C418").Value = "0196"
$ws.Range("D418").Value = "QES"
$ws.Range("E418").Value = 0.23
$ws.Range("F418").Value = 0.23
$ws.Range("G418").Value = 0.215
$ws.Range("H418").Value = 0.215
$ws.Range("I418").Value = 14167100

$ws.Range("A419").Value = 1574640000
$ws.Range("B419").Value = "2019-11-25"
$ws.Range("C419").Value = "0196"
$ws.Range("D419").Value = "QES"
$ws.Range("E419").Value = 0.215
$ws.Range("F419").Value = 0.225
$ws.Range("G419").Value = 0.21
$ws.Range("H419").Value = 0.22
$ws.Range("I419").Value = 7076400

$ws.Range("A420").Value = 1574726400
$ws.Range("B420").Value = "2019-11-26"
$ws.Range("C420").Value = "0196"
$ws.Range("D420").Value = "QES"
$ws.Range("E420").Value = 0.22
$ws.Range("F420").Value = 0.225
$ws.Range("G420").Value = 0.21
$ws.Range("H420").Value = 0.215
$ws.Range("I420").Value = 6071700

$ws.Range("A421").Value = 1574812800
$ws.Range("B421").Value = "2019-11-27"
$ws.Range("C421").Value = "0196"
$ws.Range("D421").Value = "QES"
$ws.Range("E421").Value = 0.21
$ws.Range("F421").Value = 0.225
$ws.Range("G421").Value = 0.21
$ws.Range("H421").Value = 0.22
$ws.Range("I421").Value = 10113500

$ws.Range("A422").Value = 1574899200
$ws.Range("B422").Value = "2019-11-28"
$ws.Range("C422").Value = "0196"
$ws.Range("D422").Value = "QES"
$ws.Range("E422").Value = 0.225
$ws.Range("F422").Value = 0.23
$ws.Range("G422").Value = 0.22
$ws.Range("H422").Value = 0.225
$ws.Range("I422").Value = 8201600

$ws.Range("B414:D422").Style = "Normal"